# Auto-generated edit script applying numeric updates to the Leve profit sheets
# as described by the commit diff (scheduled data refresh of currentAveragePrice /
# LevePriceNQ/HQ / LeveProfitNQ/HQ columns H:N across all 8 job sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("J17").Value = 10619.538
$ws.Range("H17").Value = 10619.538
$ws.Range("L17").Value = 31858.614
$ws.Range("K17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("N17").Value = -32194.614
$ws.Range("I18").Value = 1564
$ws.Range("K18").Value = 1564
$ws.Range("M18").Value = -1280
$ws.Range("H18").Value = 1536.25
$ws.Range("I41").Value = 3936.7334
$ws.Range("N41").Value = -6026.1665
$ws.Range("M41").Value = -3496.7334
$ws.Range("J41").Value = 5146.1665
$ws.Range("H41").Value = 4282.2856
$ws.Range("L41").Value = 5146.1665
$ws.Range("K41").Value = 3936.7334
$ws.Range("J58").Value = 4995
$ws.Range("H58").Value = 2362.2222
$ws.Range("L58").Value = 14985
$ws.Range("N58").Value = -15285
$ws.Range("N64").Value = -10495
$ws.Range("J64").Value = 9999
$ws.Range("H64").Value = 9999
$ws.Range("K64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("L64").Value = 9999
$ws.Range("J67").Value = 9999
$ws.Range("H67").Value = 9999
$ws.Range("L67").Value = 9999
$ws.Range("K67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("N67").Value = -11715
$ws.Range("N99").Value = -8358.9998
$ws.Range("M99").Value = -772.1000000000004
$ws.Range("J99").Value = 1787.6666
$ws.Range("H99").Value = 1245.0526
$ws.Range("K99").Value = 2270.1
$ws.Range("I99").Value = 756.7
$ws.Range("L99").Value = 5362.9998
$ws.Range("I101").Value = 834.5
$ws.Range("K101").Value = 2503.5
$ws.Range("M101").Value = -881.5
$ws.Range("H101").Value = 1091.8334
$ws.Range("I137").Value = 2536.7693
$ws.Range("M137").Value = -5060.3079
$ws.Range("H137").Value = 2466.7368
$ws.Range("K137").Value = 7610.3079
$ws.Range("M2").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -4629.465
$ws.Range("J32").Value = 14892.655
$ws.Range("H32").Value = 7809.56
$ws.Range("I32").Value = 4916.465
$ws.Range("L32").Value = 14892.655
$ws.Range("K32").Value = 4916.465
$ws.Range("N32").Value = -15466.655
$ws.Range("M45").Value = -5267.3335
$ws.Range("J45").Value = 4983.3335
$ws.Range("H45").Value = 5479.0835
$ws.Range("K45").Value = 5644.3335
$ws.Range("I45").Value = 5644.3335
$ws.Range("L45").Value = 4983.3335
$ws.Range("N45").Value = -5737.3335
$ws.Range("M61").Value = -3183.1538
$ws.Range("J61").Value = 4213.4
$ws.Range("H61").Value = 3750.913
$ws.Range("K61").Value = 3395.1538
$ws.Range("I61").Value = 3395.1538
$ws.Range("L61").Value = 4213.4
$ws.Range("N61").Value = -4637.4
$ws.Range("H74").Value = 2678.875
$ws.Range("K74").Value = 2419.7646
$ws.Range("I74").Value = 2419.7646
$ws.Range("M74").Value = -1545.7646
$ws.Range("M77").Value = -7730.823
$ws.Range("H77").Value = 2678.875
$ws.Range("I77").Value = 2419.7646
$ws.Range("K77").Value = 12098.823
$ws.Range("J97").Value = 1628.8
$ws.Range("H97").Value = 795.44446
$ws.Range("K97").Value = 606.0454999999999
$ws.Range("I97").Value = 606.0454999999999
$ws.Range("L97").Value = 1628.8
$ws.Range("N97").Value = -2620.8
$ws.Range("M97").Value = -110.0454999999999
$ws.Range("J132").Value = 3831.077
$ws.Range("H132").Value = 123428.96
$ws.Range("I132").Value = 177982.39
$ws.Range("L132").Value = 11493.231
$ws.Range("K132").Value = 533947.17
$ws.Range("N132").Value = -16553.231
$ws.Range("M132").Value = -531417.17
$ws.Range("H136").Value = 3750.913
$ws.Range("L136").Value = 12640.2
$ws.Range("K136").Value = 10185.4614
$ws.Range("I136").Value = 3395.1538
$ws.Range("N136").Value = -17740.2
$ws.Range("M136").Value = -7635.4614
$ws.Range("J136").Value = 4213.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M26").Value = -53041
$ws.Range("J26").Value = 104999.75
$ws.Range("H26").Value = 82856.86
$ws.Range("K26").Value = 53333
$ws.Range("I26").Value = 53333
$ws.Range("L26").Value = 104999.75
$ws.Range("N26").Value = -105583.75
$ws.Range("N28").Value = -74588
$ws.Range("J28").Value = 74000
$ws.Range("H28").Value = 74000
$ws.Range("L28").Value = 74000
$ws.Range("J107").Value = 3466.6667
$ws.Range("H107").Value = 1439.3334
$ws.Range("L107").Value = 3466.6667
$ws.Range("K107").Value = 932.5
$ws.Range("I107").Value = 932.5
$ws.Range("M107").Value = 987.5
$ws.Range("N107").Value = -7306.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L31").Value = 10581.521
$ws.Range("N31").Value = -11171.521
$ws.Range("J31").Value = 10581.521
$ws.Range("H31").Value = 9850.967000000001
$ws.Range("L34").Value = 10581.521
$ws.Range("N34").Value = -10985.521
$ws.Range("H34").Value = 9850.967000000001
$ws.Range("J34").Value = 10581.521
$ws.Range("H58").Value = 3384.8948
$ws.Range("I58").Value = 3011.2
$ws.Range("K58").Value = 3011.2
$ws.Range("M58").Value = -2808.2
$ws.Range("J119").Value = 66998.5
$ws.Range("H119").Value = 66998.5
$ws.Range("L119").Value = 66998.5
$ws.Range("N119").Value = -76674.5
$ws.Range("H132").Value = 1184.8572
$ws.Range("I132").Value = 1132.3334
$ws.Range("K132").Value = 3397.0002
$ws.Range("M132").Value = -867.0001999999999
$ws.Range("M134").Value = -355.7142000000003
$ws.Range("H134").Value = 1384.4889
$ws.Range("I134").Value = 963.5714
$ws.Range("K134").Value = 2890.7142
$ws.Range("H136").Value = 3384.8948
$ws.Range("K136").Value = 9033.599999999999
$ws.Range("I136").Value = 3011.2
$ws.Range("M136").Value = -6483.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2522.5417
$ws.Range("I132").Value = 2224.6428
$ws.Range("K132").Value = 20021.7852
$ws.Range("M132").Value = -17491.7852
$ws.Range("I137").Value = 2100
$ws.Range("M137").Value = -1200
$ws.Range("H137").Value = 2100
$ws.Range("K137").Value = 6300

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L57").Value = 34221.5
$ws.Range("N57").Value = -35861.5
$ws.Range("J57").Value = 34221.5
$ws.Range("H57").Value = 20907.166
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -31884
$ws.Range("J117").Value = 25000
$ws.Range("H117").Value = 25000
$ws.Range("H132").Value = 1759.2084
$ws.Range("I132").Value = 1541.05
$ws.Range("K132").Value = 4623.15
$ws.Range("M132").Value = -2093.15

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K40").Value = 11908828
$ws.Range("I40").Value = 11908828
$ws.Range("L40").Value = 4739.6665
$ws.Range("N40").Value = -5011.6665
$ws.Range("M40").Value = -11908692
$ws.Range("J40").Value = 4739.6665
$ws.Range("M93").Value = 587
$ws.Range("H93").Value = 1034.5883
$ws.Range("K93").Value = 661
$ws.Range("I93").Value = 661
$ws.Range("J103").Value = 31600
$ws.Range("H103").Value = 31600
$ws.Range("L103").Value = 31600
$ws.Range("N103").Value = -33944
$ws.Range("H132").Value = 439050
$ws.Range("I132").Value = 836187.5
$ws.Range("K132").Value = 2508562.5
$ws.Range("M132").Value = -2506032.5
$ws.Range("H136").Value = 5733.2188
$ws.Range("L136").Value = 21014.625
$ws.Range("K136").Value = 15928.0005
$ws.Range("I136").Value = 5309.3335
$ws.Range("N136").Value = -26114.625
$ws.Range("M136").Value = -13378.0005
$ws.Range("J136").Value = 7004.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I18").Value = 10006
$ws.Range("K18").Value = 10006
$ws.Range("M18").Value = -9833
$ws.Range("H18").Value = 16673.334
$ws.Range("M26").Value = -4707
$ws.Range("J26").Value = 12007
$ws.Range("H26").Value = 9671.333000000001
$ws.Range("K26").Value = 5000
$ws.Range("I26").Value = 5000
$ws.Range("L26").Value = 12007
$ws.Range("N26").Value = -12593
$ws.Range("J62").Value = 5607.1113
$ws.Range("H62").Value = 5607.1113
$ws.Range("L62").Value = 5607.1113
$ws.Range("N62").Value = -6855.1113
$ws.Range("H65").Value = 5607.1113
$ws.Range("L65").Value = 28035.5565
$ws.Range("N65").Value = -34275.5565
$ws.Range("J65").Value = 5607.1113
$ws.Range("J119").Value = 66122.25
$ws.Range("H119").Value = 66122.25
$ws.Range("L119").Value = 66122.25
$ws.Range("N119").Value = -75798.25
$ws.Range("J132").Value = 7356.6
$ws.Range("H132").Value = 31922.572
$ws.Range("I132").Value = 36016.9
$ws.Range("L132").Value = 22069.8
$ws.Range("K132").Value = 108050.7
$ws.Range("N132").Value = -27129.8
$ws.Range("M132").Value = -105520.7
$ws.Range("H136").Value = 57433.527
$ws.Range("L136").Value = 265609.242
$ws.Range("K136").Value = 12342.8568
$ws.Range("I136").Value = 4114.2856
$ws.Range("N136").Value = -270709.242
$ws.Range("M136").Value = -9792.856800000001
$ws.Range("J136").Value = 88536.414
